# Harmonize similar tags to be the same (Tags table on SwateTemplateMetadata sheet).
# - "phenotyping " / "study " (trailing-space variants) -> "phenotyping" / "study"
# - The free-text NCIT URL + "NCIT" pair under the "Plant" tag is replaced by a
#   proper per-tag Term Source REF / Term Accession Number layout:
#     phenotyping -> DPBO:1000224
#     Plant       -> NCIT:C14258
#     study       -> NCIT:C63536

# Note: new shared-string entries are appended in the order the cells below
# are written, so the write order here is chosen to reproduce the same
# shared-string table ordering as the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Row 13 now carries one accession number per tag column (D, E, G),
# instead of a single free-text URL crammed into E13.
$ws.Range("E13").Value = "NCIT:C14258"

# Row 12 holds the Tags themselves; harmonize the trailing-space variants.
$ws.Range("G12").Value = "study"

$ws.Range("G13").Value = "NCIT:C63536"

$ws.Range("D12").Value = "phenotyping"

$ws.Range("D13").Value = "DPBO:1000224"

# The old Term Source REF value ("NCIT") in E14 is no longer needed.
$ws.Range("E14").Value = ""

# Row 13 shrinks back to a two-line row now that the long URL is gone.
$ws.Rows.Item(13).RowHeight = 28.8

# Update the remembered selection to match where the edit was made.
$ws.Range("D16").Select()
